# "Finished WeiSSVM real data, started on HTLP"
#
# - rename the original sheet "Ark1" -> "WeiSSVM" (its data/layout is final
#   and untouched)
# - duplicate it into a second sheet "HTLP" (same BIC/AIC template/formulas)
#   and fill in the new real data that's ready so far (Summer/Fall
#   neg log-lik + df for the "cl+dfcl" and "bl+dfbl" blocks); the
#   "cl+dfbl" block's raw inputs aren't ready yet so those cells stay blank
# - leave the selection/active-sheet where the author was last working
#   (HTLP!C15)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "WeiSSVM"

# Duplicate WeiSSVM (keeps the shared formulas/number formats/row heights
# identical) right after itself, then rename the copy.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "HTLP"

# --- cl+dfcl block (rows 2-9): new real data for Summer/Fall ---
$ws2.Range("B3").Value = 878.13
$ws2.Range("C3").Value = 727.78
$ws2.Range("B4").Value = 57.94
$ws2.Range("C4").Value = 131.54

$ws2.Range("B7").Value = 1040.37
$ws2.Range("C7").Value = 821.15
$ws2.Range("B8").Value = 75.78
$ws2.Range("C8").Value = 120.88

# --- bl+dfbl block (rows 14-21): new real data for Summer/Fall ---
$ws2.Range("B15").Value = 958.24
$ws2.Range("C15").Value = 1337.51
$ws2.Range("B16").Value = 106.5
$ws2.Range("C16").Value = 65.67

$ws2.Range("B19").Value = 972.5
$ws2.Range("C19").Value = 1326.86
$ws2.Range("B20").Value = 6.5
$ws2.Range("C20").Value = 38.91

# Columns D/E (3rd/4th input series) for these two blocks aren't filled in
# yet, unlike the finished WeiSSVM sheet - clear them (the cl+dfbl block
# below them just references B/C so it's unaffected).
$ws2.Range("D3").ClearContents()
$ws2.Range("E3").ClearContents()
$ws2.Range("D4").ClearContents()
$ws2.Range("E4").ClearContents()

$ws2.Range("D7").ClearContents()
$ws2.Range("E7").ClearContents()
$ws2.Range("D8").ClearContents()
$ws2.Range("E8").ClearContents()

$ws2.Range("D15").ClearContents()
$ws2.Range("E15").ClearContents()
$ws2.Range("D16").ClearContents()
$ws2.Range("E16").ClearContents()

$ws2.Range("D19").ClearContents()
$ws2.Range("E19").ClearContents()
$ws2.Range("D20").ClearContents()
$ws2.Range("E20").ClearContents()

# Restore the cursor on the finished sheet, then leave off where work on
# the new sheet stopped.
$ws1.Range("B12").Select()
$ws2.Range("C15").Select()
